# Rename the "flock_growth"/"flock_mortality" variable labels to
# "growth"/"mortality", and add a new "perc_laying" variable for
# China and Vietnam.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing variable labels (shared strings "flock_growth" ->
# "growth", "flock_mortality" -> "mortality").
$ws.Range("B4").Value = "growth"
$ws.Range("B5").Value = "growth"
$ws.Range("B6").Value = "mortality"
$ws.Range("B7").Value = "mortality"

# Append the new perc_laying rows for China and Vietnam.
$ws.Range("A8").Value = "China"
$ws.Range("B8").Value = "perc_laying"
$ws.Range("C8").Value = 85

$ws.Range("A9").Value = "Vietnam"
$ws.Range("B9").Value = "perc_laying"
$ws.Range("C9").Value = 85

# Match the number format used by the other "value" cells below the
# header row (2 decimal places).
$ws.Range("C8:C9").NumberFormat = "0.00"

# Move the active selection to B8, matching the saved view state.
[void]$ws.Range("B8").Select()
